$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.317.05"
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = "'3.501.37"
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'590.90"
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").Value = "'133.98"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("E9").Value = '  +3.30%  '
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("D12").Value = "'4.099.06"
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").Value = "'0.0000180"
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("D15").Value = "'3.502.20"
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("D16").Value = "'64.400.43"
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = "'25.65"
$ws.Range("E17").Value = '  -6.51%  '
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("E19").Value = '  +2.44%  '
$ws.Range("E20").Value = '  -2.49%  '
$ws.Range("D21").Value = "'393.00"
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").Value = "'0.572"
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("D23").Value = "'3.641.45"
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").Value = "'74.63"
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("E26").Value = '  +0.47%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("E29").Value = '  +1.52%  '
$ws.Range("D30").Value = "'8.22"
$ws.Range("E30").Value = '  -2.67%  '
$ws.Range("E31").Value = '  -6.78%  '
$ws.Range("D32").Value = "'3.522.66"
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("E33").Value = '  +5.57%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = "'23.46"
$ws.Range("E35").Value = '  -0.44%  '
$ws.Range("D36").Value = "'5.13"
$ws.Range("E36").Value = '  -4.78%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = "'6.87"
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = "'1.55"
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("D39").Value = "'167.30"
$ws.Range("E39").Value = '  +5.53%  '
$ws.Range("D40").Value = "'0.0780"
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = '  -0.33%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = "'25.19"
$ws.Range("E43").Value = '  -5.50%  '
$ws.Range("D44").Value = "'4.39"
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E45").Value = '  +3.57%  '
$ws.Range("E46").Value = '  -3.28%  '
$ws.Range("E47").Value = '  -0.66%  '
$ws.Range("D48").Value = "'2.378.72"
$ws.Range("E48").Value = '  -4.19%  '
$ws.Range("D49").Value = "'0.892"
$ws.Range("E49").Value = '  -2.40%  '
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("D51").Value = "'21.08"
$ws.Range("E51").Value = '  -1.27%  '